$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 324.375
$ws.Range("I55").Value = 279.6
$ws.Range("J55").Value = 399
$ws.Range("K55").Value = 279.6
$ws.Range("L55").Value = 399
$ws.Range("M55").Value = -65.60000000000002
$ws.Range("N55").Value = -827

$ws.Range("H70").Value = 111223920
$ws.Range("J70").Value = 125126750
$ws.Range("L70").Value = 375380250
$ws.Range("N70").Value = -375380790

$ws.Range("H73").Value = 111223920
$ws.Range("J73").Value = 125126750
$ws.Range("L73").Value = 375380250
$ws.Range("N73").Value = -375382122

$ws.Range("H106").Value = 5721.0527
$ws.Range("I106").Value = 6400
$ws.Range("K106").Value = 6400
$ws.Range("M106").Value = -5769

$ws.Range("H113").Value = 4699.75
$ws.Range("I113").Value = 5000
$ws.Range("J113").Value = 4599.6665
$ws.Range("K113").Value = 5000
$ws.Range("L113").Value = 4599.6665
$ws.Range("M113").Value = -1746
$ws.Range("N113").Value = -11107.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5070.1177
$ws.Range("I2").Value = 5574.6665
$ws.Range("K2").Value = 5574.6665
$ws.Range("M2").Value = -5461.6665

$ws.Range("H74").Value = 2562.18
$ws.Range("I74").Value = 2536
$ws.Range("K74").Value = 2536
$ws.Range("M74").Value = -1662

$ws.Range("H77").Value = 2562.18
$ws.Range("I77").Value = 2536
$ws.Range("K77").Value = 12680
$ws.Range("M77").Value = -8312

$ws.Range("H97").Value = 604.6896400000001
$ws.Range("I97").Value = 595.52
$ws.Range("J97").Value = 662
$ws.Range("K97").Value = 595.52
$ws.Range("L97").Value = 662
$ws.Range("M97").Value = -99.51999999999998
$ws.Range("N97").Value = -1654

$ws.Range("H110").Value = 2570.75
$ws.Range("I110").Value = 2219.5
$ws.Range("K110").Value = 2219.5
$ws.Range("M110").Value = -174.5

$ws.Range("H116").Value = 5070.1177
$ws.Range("I116").Value = 5574.6665
$ws.Range("K116").Value = 5574.6665
$ws.Range("M116").Value = -3280.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5070.1177
$ws.Range("I3").Value = 5574.6665
$ws.Range("K3").Value = 5574.6665
$ws.Range("M3").Value = -5460.6665

$ws.Range("H99").Value = 3111.5789
$ws.Range("J99").Value = 3948.75
$ws.Range("L99").Value = 3948.75
$ws.Range("N99").Value = -6944.75

$ws.Range("H105").Value = 2556.52
$ws.Range("I105").Value = 2037.0526
$ws.Range("K105").Value = 2037.0526
$ws.Range("M105").Value = -290.0526

$ws.Range("H134").Value = 2470.5925
$ws.Range("I134").Value = 2238.7114
$ws.Range("K134").Value = 6716.1342
$ws.Range("M134").Value = -4181.1342

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 43415.145
$ws.Range("I16").Value = 32358.8
$ws.Range("J16").Value = 71056
$ws.Range("K16").Value = 32358.8
$ws.Range("L16").Value = 71056
$ws.Range("M16").Value = -32071.8
$ws.Range("N16").Value = -71630

$ws.Range("H58").Value = 2503.5557
$ws.Range("I58").Value = 3286.3333
$ws.Range("K58").Value = 3286.3333
$ws.Range("M58").Value = -3083.3333

$ws.Range("H105").Value = 2835.4
$ws.Range("I105").Value = 2966
$ws.Range("K105").Value = 2966
$ws.Range("M105").Value = -1219

$ws.Range("H113").Value = 43415.145
$ws.Range("I113").Value = 32358.8
$ws.Range("J113").Value = 71056
$ws.Range("K113").Value = 32358.8
$ws.Range("L113").Value = 71056
$ws.Range("M113").Value = -30188.8
$ws.Range("N113").Value = -75396

$ws.Range("H136").Value = 2503.5557
$ws.Range("I136").Value = 3286.3333
$ws.Range("K136").Value = 9858.999899999999
$ws.Range("M136").Value = -7308.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 19066.166
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 19066.166
$ws.Range("K76").Value = 0
$ws.Range("L76").ClearContents()
$ws.Range("M76").Value = 57198.49800000001
$ws.Range("N76").Value = -57964.49800000001

$ws.Range("H79").Value = 19066.166
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 19066.166
$ws.Range("K79").Value = 0
$ws.Range("L79").ClearContents()
$ws.Range("M79").Value = 57198.49800000001
$ws.Range("N79").Value = -59850.49800000001

$ws.Range("H109").Value = 2225.6365
$ws.Range("I109").Value = 1468.4706
$ws.Range("J109").Value = 4800
$ws.Range("K109").Value = 4405.4118
$ws.Range("L109").Value = 14400
$ws.Range("M109").Value = -3365.4118
$ws.Range("N109").Value = -16480

$ws.Range("H110").Value = 2513.5
$ws.Range("I110").Value = 2513.5
$ws.Range("K110").Value = 7540.5
$ws.Range("M110").Value = -3450.5

$ws.Range("H111").Value = 762
$ws.Range("I111").Value = 762
$ws.Range("K111").Value = 2286
$ws.Range("M111").Value = 781

$ws.Range("H112").Value = 6440.3887
$ws.Range("J112").Value = 6636.364
$ws.Range("L112").Value = 19909.092
$ws.Range("N112").Value = -22125.092

$ws.Range("H116").Value = 1645.8
$ws.Range("I116").Value = 1670.25
$ws.Range("K116").Value = 5010.75
$ws.Range("M116").Value = -1568.75

$ws.Range("H118").Value = 5024.9375
$ws.Range("I118").Value = 1301
$ws.Range("K118").Value = 3903
$ws.Range("M118").Value = -2660

$ws.Range("H119").Value = 3285.4443
$ws.Range("I119").Value = 3285.4443
$ws.Range("K119").Value = 9856.332900000001
$ws.Range("M119").Value = -5018.332900000001

$ws.Range("H120").Value = 10826.333
$ws.Range("I120").Value = 6991.8
$ws.Range("K120").Value = 20975.4
$ws.Range("M120").Value = -16137.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6825.533
$ws.Range("I113").Value = 7630.4
$ws.Range("K113").Value = 7630.4
$ws.Range("M113").Value = -5460.4

$ws.Range("H122").Value = 2376.2222
$ws.Range("I122").Value = 2376.2222
$ws.Range("K122").Value = 7128.6666
$ws.Range("M122").Value = -4678.6666

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("N141").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2210.5945
$ws.Range("I136").Value = 2076.125
$ws.Range("K136").Value = 6228.375
$ws.Range("M136").Value = -3678.375

$ws.Range("H137").Value = 89451.73
$ws.Range("J137").Value = 89451.73
$ws.Range("L137").Value = 89451.73
$ws.Range("N137").Value = -99651.73

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 45000
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

$ws.Range("H100").Value = 2894.9167
$ws.Range("I100").Value = 2921.7144
$ws.Range("J100").Value = 2857.4
$ws.Range("K100").Value = 5843.4288
$ws.Range("L100").Value = 5714.8
$ws.Range("M100").Value = -5302.4288
$ws.Range("N100").Value = -6796.8

$ws.Range("H114").Value = 55741
$ws.Range("J114").Value = 55741
$ws.Range("L114").Value = 55741
$ws.Range("N114").Value = -64419

$ws.Range("H122").Value = 2629.8
$ws.Range("I122").Value = 2059.5715
$ws.Range("K122").Value = 6178.7145
$ws.Range("M122").Value = -3728.7145

$ws.Range("H132").Value = 2324.5715
$ws.Range("I132").Value = 2336.6296
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 7009.888800000001
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -4479.888800000001
$ws.Range("N132").Value = -11057

$ws.Range("H136").Value = 2491.7354
$ws.Range("I136").Value = 2230.3704
$ws.Range("K136").Value = 6691.111199999999
$ws.Range("M136").Value = -4141.111199999999
